# Records trade #2 (base_strategy, UP, OPEN) as a new row 3 on both the
# "All Trades" and "base_strategy" worksheets - mirrors how trade #1 is
# already laid out in row 2 of each sheet.
#
# Notes on technique:
#  - Plain `.Value = "2026-02-16"` / `"22:57:24"`-style assignments get
#    smart-typed by Excel into date/time serials, which isn't what the
#    source data is (it's stored as literal text). Round-tripping the
#    string through a text formula (`="..."`) and then Copy +
#    PasteSpecial(xlPasteValues) "bakes" the formula result back down to a
#    plain literal cell without Excel re-inferring a date/number type and
#    without leaving a custom number-format style behind.
#  - Plain numbers are written directly since they don't have this
#    ambiguity.

$xlPasteValues = -4163

function Add-TradeRow2 {
    param($ws)

    $ws.Range("A3").Value = 2

    $ws.Range("B3").Formula = '="2026-02-16"'
    $ws.Range("B3").Copy()
    $ws.Range("B3").PasteSpecial($xlPasteValues)

    $ws.Range("C3").Formula = '="22:57:24"'
    $ws.Range("C3").Copy()
    $ws.Range("C3").PasteSpecial($xlPasteValues)

    $ws.Range("D3").Formula = '="base_strategy"'
    $ws.Range("D3").Copy()
    $ws.Range("D3").PasteSpecial($xlPasteValues)

    $ws.Range("E3").Formula = '="UP"'
    $ws.Range("E3").Copy()
    $ws.Range("E3").PasteSpecial($xlPasteValues)

    $ws.Range("F3").Value = 0.5

    $ws.Range("G3").Formula = '=""'
    $ws.Range("G3").Copy()
    $ws.Range("G3").PasteSpecial($xlPasteValues)

    $ws.Range("H3").Formula = '="OPEN"'
    $ws.Range("H3").Copy()
    $ws.Range("H3").PasteSpecial($xlPasteValues)

    $ws.Range("I3").Value = 0
    $ws.Range("J3").Value = 0
    $ws.Range("K3").Value = 100
    $ws.Range("L3").Value = 0
    $ws.Range("M3").Value = 0
    $ws.Range("N3").Value = 0.6

    $ws.Range("O3").Formula = '="Normal spread capture: 19600 bps"'
    $ws.Range("O3").Copy()
    $ws.Range("O3").PasteSpecial($xlPasteValues)

    $ws.Range("P3").Formula = '=""'
    $ws.Range("P3").Copy()
    $ws.Range("P3").PasteSpecial($xlPasteValues)

    $ws.Range("Q3").Value = 0
}

$wb = $excel.ActiveWorkbook

$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow2 $wsAllTrades

$wsBaseStrategy = $wb.Worksheets.Item("base_strategy")
Add-TradeRow2 $wsBaseStrategy
